$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 169, which pushes the existing rows 169-174 down to 170-175
# (matching content, styles such as the date format on column D are inherited
# automatically from the row being pushed down).
$ws.Rows.Item(169).Insert()

# Populate the newly inserted row 169 with the new weekly price record.
$ws.Cells.Item(169, 1).Value = 8
$ws.Cells.Item(169, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(169, 3).Value = "Coquimbo"
$ws.Cells.Item(169, 4).Value = 45041
$ws.Cells.Item(169, 5).Value = 4
$ws.Cells.Item(169, 6).Value = 100112052
$ws.Cells.Item(169, 7).Value = "Albahaca"
$ws.Cells.Item(169, 8).Value = "Sin especificar"
$ws.Cells.Item(169, 9).Value = "Primera"
$ws.Cells.Item(169, 10).Value = 600
$ws.Cells.Item(169, 11).Value = 2800
$ws.Cells.Item(169, 12).Value = 3000
$ws.Cells.Item(169, 13).Value = 2900
$ws.Cells.Item(169, 14).Value = "`$/paquete"
$ws.Cells.Item(169, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(169, 16).Value = 2900
$ws.Cells.Item(169, 17).Value = 1
$ws.Cells.Item(169, 18).Value = "Hortaliza"

# Ensure column D keeps the date number format used by the rest of the column.
$ws.Cells.Item(169, 4).NumberFormat = $ws.Cells.Item(170, 4).NumberFormat
